$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric rows 2-7 (column A = Meta, column B = Venda)
$ws.Range("A2").Value = 3000.6
$ws.Range("B2").Value = 6000

$ws.Range("A3").Value = 9800
$ws.Range("B3").Value = 12500

$ws.Range("A4").Value = 9000
$ws.Range("B4").Value = 1000

$ws.Range("A5").Value = 6000
$ws.Range("B5").Value = 15000

$ws.Range("A6").Value = 6700
$ws.Range("B6").Value = 0

$ws.Range("A7").Value = 0.6
$ws.Range("B7").Value = 0

# Row 8 holds text-formatted numbers (quote-prefixed so they stay text, not numbers)
$ws.Range("A8").Value = "'1.2"
$ws.Range("B8").Value = "'0.0"
